$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text entry (avoid Excel re-parsing numeric-looking strings like "1.001" as numbers)
# then restore the default "Normal" style so no style/format change is recorded.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '25.787.30'
$ws.Range("E2").Value = '  -3.80%  '

$ws.Range("D3").Value = '1.820.28'
$ws.Range("E3").Value = '  -2.81%  '

$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Value = '278.51'
$ws.Range("E5").Value = '  -7.55%  '

$ws.Range("E6").Value = '  -0.15%  '

$ws.Range("D7").Value = '0.5093'
$ws.Range("E7").Value = '  -4.45%  '

$ws.Range("E8").Value = '  -5.32%  '

$ws.Range("D9").Value = '44.55'
$ws.Range("E9").Value = '  -2.18%  '

$ws.Range("D10").Value = '0.06668'
$ws.Range("E10").Value = '  -7.19%  '

$ws.Range("D11").Value = '20.05'
$ws.Range("E11").Value = '  -7.27%  '

$ws.Range("D12").Value = '0.8277'
$ws.Range("E12").Value = '  -6.90%  '

$ws.Range("D13").Value = '0.07903'
$ws.Range("E13").Value = '  -3.18%  '

$ws.Range("D14").Value = '1.821.83'
$ws.Range("E14").Value = '  -2.92%  '

$ws.Range("D15").Value = '5.082'
$ws.Range("E15").Value = '  -4.01%  '

$ws.Range("D16").Value = '88.12'
$ws.Range("E16").Value = '  -5.17%  '

$ws.Range("E17").Value = '  -0.10%  '

$ws.Range("E18").Value = '  -4.83%  '

$ws.Range("D19").Value = '0.000008034'
$ws.Range("E19").Value = '  -5.66%  '

$ws.Range("E20").Value = '  -0.13%  '

$ws.Range("D21").Value = '25.824.77'
$ws.Range("E21").Value = '  -3.82%  '

$ws.Range("D22").Value = '4.758'
$ws.Range("E22").Value = '  -4.53%  '

$ws.Range("D23").Value = '10.01'
$ws.Range("E23").Value = '  -5.90%  '

$ws.Range("D24").Value = '6.119'
$ws.Range("E24").Value = '  -4.34%  '

$ws.Range("D25").Value = '2.237'
$ws.Range("E25").Value = '  -2.29%  '

$ws.Range("D26").Value = '142.29'
$ws.Range("E26").Value = '  -2.86%  '

$ws.Range("D27").Value = '1.671'
$ws.Range("E27").Value = '  -4.06%  '

$ws.Range("D28").Value = '17.15'
$ws.Range("E28").Value = '  -4.84%  '

$ws.Range("E29").Value = '  -4.00%  '

$ws.Range("D30").Value = '4.324'
$ws.Range("E30").Value = '  -8.26%  '

$ws.Range("D31").Value = '4.243'
$ws.Range("E31").Value = '  -8.02%  '

$ws.Range("D32").Value = '0.08750'
$ws.Range("E32").Value = '  -4.13%  '

$ws.Range("D33").Value = '0.04915'
$ws.Range("E33").Value = '  -1.91%  '

$ws.Range("D34").Value = '0.7306'
$ws.Range("E34").Value = '  -9.93%  '

$ws.Range("E35").Value = '  -2.64%  '

$ws.Range("D36").Value = '2.883'
$ws.Range("E36").Value = '  -2.77%  '

$ws.Range("D37").Value = '3.159'
$ws.Range("E37").Value = '  -1.36%  '

$ws.Range("E38").Value = '  -0.26%  '

$ws.Range("D39").Value = '2.364'
$ws.Range("E39").Value = '  -11.13%  '

$ws.Range("D40").Value = '0.01861'

$ws.Range("D41").Value = '0.5160'
$ws.Range("E41").Value = '  -14.88%  '

$ws.Range("D42").Value = '0.9709'
$ws.Range("E42").Value = '  -9.13%  '

$ws.Range("D43").Value = '114.06'
$ws.Range("E43").Value = '  -0.98%  '

$ws.Range("D44").Value = '6.249'
$ws.Range("E44").Value = '  -4.62%  '

$ws.Range("D45").Value = '8.049'
$ws.Range("E45").Value = '  -8.92%  '

$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  -0.17%  '

$ws.Range("D47").Value = '0.4555'
$ws.Range("E47").Value = '  -11.65%  '

$ws.Range("D48").Value = '0.1372'
$ws.Range("E48").Value = '  -8.24%  '

$ws.Range("D49").Value = '36.63'
$ws.Range("E49").Value = '  -2.43%  '

$ws.Range("D50").Value = '9.259'
$ws.Range("E50").Value = '  -6.71%  '

$ws.Range("D51").Value = '1.501'
$ws.Range("E51").Value = '  -8.49%  '

$ws.Range("D2:E51").Style = "Normal"
